$d = $word.ActiveDocument

$replacements = @(
    @{old="310×6="; new="106×9="},
    @{old="311×6="; new="848×9="},
    @{old="176×8="; new="901×7="},
    @{old="288×7="; new="955×8="},
    @{old="222×5="; new="515×9="},
    @{old="257×2="; new="734×7="},
    @{old="716×7="; new="872×4="},
    @{old="740×5="; new="428×3="},
    @{old="873×3="; new="505×6="},
    @{old="392×7="; new="428×2="},
    @{old="493×9="; new="423×8="},
    @{old="150×9="; new="985×8="},
    @{old="509×6="; new="480×3="},
    @{old="229×2="; new="286×2="},
    @{old="598×4="; new="876×3="},
    @{old="772×4="; new="571×4="},
    @{old="288×5="; new="534×4="},
    @{old="845×5="; new="539×7="},
    @{old="266×4="; new="993×2="},
    @{old="538×8="; new="465×4="},
    @{old="410×4="; new="418×7="},
    @{old="894×7="; new="987×6="},
    @{old="227×2="; new="296×8="},
    @{old="676×6="; new="103×3="},
    @{old="141×8="; new="909×6="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
